$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.718.67'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '2.053.10'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.46'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -1.59%  '
$ws.Range("E10").Value = '  +3.53%  '
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").Value = '2.354.03'
$ws.Range("E12").Value = '  +0.99%  '
$ws.Range("E13").Value = '  -0.67%  '
$ws.Range("E14").Value = '  +1.52%  '
$ws.Range("E15").Value = '  +6.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.763'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("D17").Value = '2.055.90'
$ws.Range("E17").Value = '  +1.31%  '
$ws.Range("D18").Value = '37.732.43'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.95'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.78%  '
$ws.Range("E21").Value = '  +1.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '222.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.87%  '
$ws.Range("E24").Value = '  +0.50%  '
$ws.Range("E25").Value = '  +2.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.29'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("E29").Value = '  -0.63%  '
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.119'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.39%  '
$ws.Range("E32").Value = '  +9.19%  '
$ws.Range("E33").Value = '  -0.96%  '
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("E35").Value = '  +0.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.47'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.36'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.43'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.07%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.38'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.01%  '
$ws.Range("D41").Value = '1.526.77'
$ws.Range("E41").Value = '  +0.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.31'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.90%  '
$ws.Range("E43").Value = '  -1.65%  '
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("E46").Value = '  -2.45%  '
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.01'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.96'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("E50").Value = '  +0.52%  '
$ws.Range("D51").Value = '2.244.76'
$ws.Range("E51").Value = '  +1.16%  '
